# Updates the cryptos worksheet with refreshed price/volume data
# (and restores the correct Uniswap / WrappedBTC row order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell updates: Row, Column, NewValue
$updates = @(
    [PSCustomObject]@{ Cell = "D2"; Value = '70.513.47' },
    [PSCustomObject]@{ Cell = "E2"; Value = '  -0.61%  ' },
    [PSCustomObject]@{ Cell = "D3"; Value = '3.552.41' },
    [PSCustomObject]@{ Cell = "E3"; Value = '  -0.95%  ' },
    [PSCustomObject]@{ Cell = "E4"; Value = '  -0.20%  ' },
    [PSCustomObject]@{ Cell = "D5"; Value = '616.79' },
    [PSCustomObject]@{ Cell = "E5"; Value = '  +5.09%  ' },
    [PSCustomObject]@{ Cell = "D6"; Value = '186.65' },
    [PSCustomObject]@{ Cell = "E6"; Value = '  +0.86%  ' },
    [PSCustomObject]@{ Cell = "D7"; Value = '0.629' },
    [PSCustomObject]@{ Cell = "E7"; Value = '  +1.20%  ' },
    [PSCustomObject]@{ Cell = "E8"; Value = '  -0.10%  ' },
    [PSCustomObject]@{ Cell = "D9"; Value = '0.218' },
    [PSCustomObject]@{ Cell = "E9"; Value = '  +0.50%  ' },
    [PSCustomObject]@{ Cell = "D10"; Value = '0.661' },
    [PSCustomObject]@{ Cell = "E10"; Value = '  +1.74%  ' },
    [PSCustomObject]@{ Cell = "D11"; Value = '53.78' },
    [PSCustomObject]@{ Cell = "E11"; Value = '  -0.80%  ' },
    [PSCustomObject]@{ Cell = "D12"; Value = '0.0000309' },
    [PSCustomObject]@{ Cell = "E12"; Value = '  -3.98%  ' },
    [PSCustomObject]@{ Cell = "D13"; Value = '9.65' },
    [PSCustomObject]@{ Cell = "E13"; Value = '  +1.38%  ' },
    [PSCustomObject]@{ Cell = "D14"; Value = '4.109.38' },
    [PSCustomObject]@{ Cell = "E14"; Value = '  -1.07%  ' },
    [PSCustomObject]@{ Cell = "D15"; Value = '620.38' },
    [PSCustomObject]@{ Cell = "E15"; Value = '  +9.40%  ' },
    [PSCustomObject]@{ Cell = "B16"; Value = 'Uniswap' },
    [PSCustomObject]@{ Cell = "C16"; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' },
    [PSCustomObject]@{ Cell = "D16"; Value = '12.87' },
    [PSCustomObject]@{ Cell = "E16"; Value = '  +3.74%  ' },
    [PSCustomObject]@{ Cell = "B17"; Value = 'WrappedBTC' },
    [PSCustomObject]@{ Cell = "C17"; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    [PSCustomObject]@{ Cell = "D17"; Value = '70.408.48' },
    [PSCustomObject]@{ Cell = "E17"; Value = '  -0.68%  ' },
    [PSCustomObject]@{ Cell = "D18"; Value = '19.11' },
    [PSCustomObject]@{ Cell = "E18"; Value = '  -0.98%  ' },
    [PSCustomObject]@{ Cell = "D19"; Value = '3.551.73' },
    [PSCustomObject]@{ Cell = "E19"; Value = '  -1.60%  ' },
    [PSCustomObject]@{ Cell = "E20"; Value = '  -0.02%  ' },
    [PSCustomObject]@{ Cell = "D21"; Value = '1.00' },
    [PSCustomObject]@{ Cell = "E21"; Value = '  -1.15%  ' },
    [PSCustomObject]@{ Cell = "D22"; Value = '17.63' },
    [PSCustomObject]@{ Cell = "E22"; Value = '  +0.00%  ' },
    [PSCustomObject]@{ Cell = "D23"; Value = '104.65' },
    [PSCustomObject]@{ Cell = "E23"; Value = '  +10.46%  ' },
    [PSCustomObject]@{ Cell = "D24"; Value = '4.73' },
    [PSCustomObject]@{ Cell = "E24"; Value = '  +2.60%  ' },
    [PSCustomObject]@{ Cell = "D25"; Value = '5.09' },
    [PSCustomObject]@{ Cell = "E25"; Value = '  +2.02%  ' },
    [PSCustomObject]@{ Cell = "D26"; Value = '3.04' },
    [PSCustomObject]@{ Cell = "E26"; Value = '  +4.06%  ' },
    [PSCustomObject]@{ Cell = "D27"; Value = '11.04' },
    [PSCustomObject]@{ Cell = "E27"; Value = '  -1.35%  ' },
    [PSCustomObject]@{ Cell = "D28"; Value = '9.86' },
    [PSCustomObject]@{ Cell = "E28"; Value = '  +8.37%  ' },
    [PSCustomObject]@{ Cell = "D29"; Value = '34.14' },
    [PSCustomObject]@{ Cell = "E29"; Value = '  +5.60%  ' },
    [PSCustomObject]@{ Cell = "D30"; Value = '7.11' },
    [PSCustomObject]@{ Cell = "E30"; Value = '  -1.86%  ' },
    [PSCustomObject]@{ Cell = "D31"; Value = '12.47' },
    [PSCustomObject]@{ Cell = "E31"; Value = '  +1.49%  ' },
    [PSCustomObject]@{ Cell = "E32"; Value = '  +1.49%  ' },
    [PSCustomObject]@{ Cell = "D33"; Value = '64.24' },
    [PSCustomObject]@{ Cell = "E33"; Value = '  +0.07%  ' },
    [PSCustomObject]@{ Cell = "D34"; Value = '3.63' },
    [PSCustomObject]@{ Cell = "E34"; Value = '  +16.54%  ' },
    [PSCustomObject]@{ Cell = "D35"; Value = '3.22' },
    [PSCustomObject]@{ Cell = "E35"; Value = '  -1.68%  ' },
    [PSCustomObject]@{ Cell = "D36"; Value = '533.72' },
    [PSCustomObject]@{ Cell = "E36"; Value = '  -2.94%  ' },
    [PSCustomObject]@{ Cell = "D37"; Value = '0.403' },
    [PSCustomObject]@{ Cell = "E37"; Value = '  -2.44%  ' },
    [PSCustomObject]@{ Cell = "E38"; Value = '  +0.10%  ' },
    [PSCustomObject]@{ Cell = "D39"; Value = '37.42' },
    [PSCustomObject]@{ Cell = "E39"; Value = '  -0.27%  ' },
    [PSCustomObject]@{ Cell = "D40"; Value = '3.62' },
    [PSCustomObject]@{ Cell = "E40"; Value = '  +6.63%  ' },
    [PSCustomObject]@{ Cell = "D41"; Value = '0.0₃0783' },
    [PSCustomObject]@{ Cell = "E41"; Value = '  -2.68%  ' },
    [PSCustomObject]@{ Cell = "D42"; Value = '3.547.35' },
    [PSCustomObject]@{ Cell = "E42"; Value = '  +2.39%  ' },
    [PSCustomObject]@{ Cell = "D43"; Value = '0.138' },
    [PSCustomObject]@{ Cell = "E43"; Value = '  +0.91%  ' },
    [PSCustomObject]@{ Cell = "D44"; Value = '0.0468' },
    [PSCustomObject]@{ Cell = "E44"; Value = '  +5.57%  ' },
    [PSCustomObject]@{ Cell = "E45"; Value = '  -0.17%  ' },
    [PSCustomObject]@{ Cell = "E46"; Value = '  +4.07%  ' },
    [PSCustomObject]@{ Cell = "E47"; Value = '  -3.51%  ' },
    [PSCustomObject]@{ Cell = "D48"; Value = '9.09' },
    [PSCustomObject]@{ Cell = "E48"; Value = '  -2.51%  ' },
    [PSCustomObject]@{ Cell = "E49"; Value = '  +0.02%  ' },
    [PSCustomObject]@{ Cell = "D50"; Value = '1.41' },
    [PSCustomObject]@{ Cell = "E50"; Value = '  -1.32%  ' },
    [PSCustomObject]@{ Cell = "D51"; Value = '134.29' },
    [PSCustomObject]@{ Cell = "E51"; Value = '  -0.32%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell -match "^D") {
        # Column D holds price text that can look numeric (e.g. "1.00", "0.0000309").
        # Force a text number-format so Excel keeps the exact literal string
        # instead of silently converting it to a float, then restore the
        # default "Normal" style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
